# Auto-generated Excel COM-interop script applying the scheduled-runner price updates
# to the per-profession Leve profit tables (Sheets/Masamune_Profits.xlsx diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 22500
$ws.Range("J13").Value = 22500
$ws.Range("L13").Value = 22500
$ws.Range("N13").Value = -22838
$ws.Range("H41").Value = 617.4167
$ws.Range("I41").Value = 94.75
$ws.Range("J41").Value = 878.75
$ws.Range("K41").Value = 94.75
$ws.Range("L41").Value = 878.75
$ws.Range("M41").Value = 345.25
$ws.Range("N41").Value = -1758.75
$ws.Range("H53").Value = 207.33333
$ws.Range("I53").Value = 88.75
$ws.Range("J53").Value = 280.30768
$ws.Range("K53").Value = 88.75
$ws.Range("L53").Value = 280.30768
$ws.Range("M53").Value = 548.25
$ws.Range("N53").Value = -1554.30768
$ws.Range("H54").Value = 54333.332
$ws.Range("J54").Value = 54000
$ws.Range("L54").Value = 54000
$ws.Range("N54").Value = -54972
$ws.Range("H58").Value = 2492.4
$ws.Range("I58").Value = 3151.1667
$ws.Range("J58").Value = 1504.25
$ws.Range("K58").Value = 9453.500100000001
$ws.Range("L58").Value = 4512.75
$ws.Range("M58").Value = -9303.500100000001
$ws.Range("N58").Value = -4812.75
$ws.Range("H86").Value = 7146182
$ws.Range("I86").Value = 9094349
$ws.Range("J86").Value = 2901.3333
$ws.Range("K86").Value = 9094349
$ws.Range("L86").Value = 2901.3333
$ws.Range("M86").Value = -9093226
$ws.Range("N86").Value = -5147.3333
$ws.Range("H89").Value = 7146182
$ws.Range("I89").Value = 9094349
$ws.Range("J89").Value = 2901.3333
$ws.Range("K89").Value = 45471745
$ws.Range("L89").Value = 14506.6665
$ws.Range("M89").Value = -45466129
$ws.Range("N89").Value = -25738.6665
$ws.Range("H105").Value = 36415.75
$ws.Range("J105").Value = 36415.75
$ws.Range("L105").Value = 36415.75
$ws.Range("N105").Value = -43403.75
$ws.Range("H129").Value = 2151.4194
$ws.Range("J129").Value = 2139.652
$ws.Range("L129").Value = 6418.956
$ws.Range("N129").Value = -16418.956
$ws.Range("H135").Value = 16668117
$ws.Range("I135").Value = 1394.7368
$ws.Range("J135").Value = 45456092
$ws.Range("K135").Value = 12552.6312
$ws.Range("L135").Value = 409104828
$ws.Range("M135").Value = -10017.6312
$ws.Range("N135").Value = -409109898
$ws.Range("H141").Value = 3019.84
$ws.Range("I141").Value = 1934.05
$ws.Range("J141").Value = 7363
$ws.Range("K141").Value = 5802.15
$ws.Range("L141").Value = 22089
$ws.Range("M141").Value = -622.1499999999996
$ws.Range("N141").Value = -32449

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30233.654
$ws.Range("I32").Value = 30516.74
$ws.Range("J32").Value = 27402.8
$ws.Range("K32").Value = 30516.74
$ws.Range("L32").Value = 27402.8
$ws.Range("M32").Value = -30229.74
$ws.Range("N32").Value = -27976.8
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H60").Value = 33990
$ws.Range("J60").Value = 33990
$ws.Range("L60").Value = 33990
$ws.Range("N60").Value = -35456
$ws.Range("H80").Value = 49124
$ws.Range("J80").Value = 49124
$ws.Range("L80").Value = 49124
$ws.Range("N80").Value = -51120
$ws.Range("H83").Value = 49124
$ws.Range("J83").Value = 49124
$ws.Range("L83").Value = 147372
$ws.Range("N83").Value = -157356
$ws.Range("H103").Value = 38338.25
$ws.Range("J103").Value = 38338.25
$ws.Range("L103").Value = 38338.25
$ws.Range("N103").Value = -40682.25
$ws.Range("H106").Value = 46496
$ws.Range("J106").Value = 46496
$ws.Range("L106").Value = 46496
$ws.Range("N106").Value = -49020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 46311.668
$ws.Range("J100").Value = 46311.668
$ws.Range("L100").Value = 46311.668
$ws.Range("N100").Value = -48475.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 43198.43
$ws.Range("J41").Value = 43198.43
$ws.Range("L41").Value = 43198.43
$ws.Range("N41").Value = -44054.43
$ws.Range("H82").Value = 37661.332
$ws.Range("J82").Value = 37661.332
$ws.Range("L82").Value = 37661.332
$ws.Range("N82").Value = -38383.332
$ws.Range("H85").Value = 37661.332
$ws.Range("J85").Value = 37661.332
$ws.Range("L85").Value = 37661.332
$ws.Range("N85").Value = -40157.332
$ws.Range("H92").Value = 37998.625
$ws.Range("J92").Value = 37998.625
$ws.Range("L92").Value = 37998.625
$ws.Range("N92").Value = -42990.625
$ws.Range("H99").Value = 1616.7084
$ws.Range("I99").Value = 1651.7
$ws.Range("J99").Value = 1591.7142
$ws.Range("K99").Value = 1651.7
$ws.Range("L99").Value = 1591.7142
$ws.Range("M99").Value = -153.7
$ws.Range("N99").Value = -4587.7142
$ws.Range("H107").Value = 664.3043
$ws.Range("I107").Value = 655.4
$ws.Range("J107").Value = 671.1539
$ws.Range("K107").Value = 655.4
$ws.Range("L107").Value = 671.1539
$ws.Range("M107").Value = 1264.6
$ws.Range("N107").Value = -4511.1539
$ws.Range("H126").Value = 1616.7084
$ws.Range("I126").Value = 1651.7
$ws.Range("J126").Value = 1591.7142
$ws.Range("K126").Value = 4955.1
$ws.Range("L126").Value = 4775.142599999999
$ws.Range("M126").Value = -2485.1
$ws.Range("N126").Value = -9715.142599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 68930.766
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = 74674.164
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = 74674.164
$ws.Range("M22").Value = 519
$ws.Range("N22").Value = -75732.164
$ws.Range("H80").Value = 138746.08
$ws.Range("I80").Value = 253740.75
$ws.Range("J80").Value = 3458.2354
$ws.Range("K80").Value = 253740.75
$ws.Range("L80").Value = 3458.2354
$ws.Range("M80").Value = -252742.75
$ws.Range("N80").Value = -5454.2354
$ws.Range("H83").Value = 138746.08
$ws.Range("I83").Value = 253740.75
$ws.Range("J83").Value = 3458.2354
$ws.Range("K83").Value = 1268703.75
$ws.Range("L83").Value = 17291.177
$ws.Range("M83").Value = -1263711.75
$ws.Range("N83").Value = -27275.177
$ws.Range("H104").Value = 37340.4
$ws.Range("J104").Value = 37340.4
$ws.Range("L104").Value = 37340.4
$ws.Range("N104").Value = -44328.4
$ws.Range("H105").Value = 45104.668
$ws.Range("J105").Value = 45104.668
$ws.Range("L105").Value = 45104.668
$ws.Range("N105").Value = -52092.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002
$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008
$ws.Range("H92").Value = 38374
$ws.Range("J92").Value = 38374
$ws.Range("L92").Value = 38374
$ws.Range("N92").Value = -43366

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 36499.5
$ws.Range("J98").Value = 36499.5
$ws.Range("L98").Value = 36499.5
$ws.Range("N98").Value = -42489.5
$ws.Range("H109").Value = 32373
$ws.Range("J109").Value = 32373
$ws.Range("L109").Value = 32373
$ws.Range("N109").Value = -35147
$ws.Range("H135").Value = 19271.412
$ws.Range("J135").Value = 19271.412
$ws.Range("L135").Value = 19271.412
$ws.Range("N135").Value = -29411.412
